$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from H1 onto the new
# header cells before setting their text, matching the style used by the
# other header cells in row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @(
    @(2,6,7),
    @(3,8,8),
    @(4,8,8),
    @(5,5,6),
    @(6,7,7),
    @(7,5,5),
    @(8,6,6),
    @(9,6,6),
    @(10,7,7),
    @(11,4,4),
    @(12,6,6),
    @(13,2,2),
    @(14,4,5),
    @(15,4,4),
    @(16,6,6),
    @(17,2,2),
    @(18,5,6),
    @(19,8,9),
    @(20,6,6),
    @(21,7,8),
    @(22,6,7),
    @(23,4,5),
    @(24,3,5),
    @(25,9,9),
    @(26,9,9),
    @(27,8,9),
    @(28,6,6),
    @(29,8,8),
    @(30,7,7),
    @(31,7,7),
    @(32,4,5),
    @(33,8,8),
    @(34,5,5),
    @(35,6,6),
    @(36,8,8),
    @(37,5,5),
    @(38,9,9),
    @(39,7,7),
    @(40,8,8),
    @(41,7,7),
    @(42,10,10),
    @(43,7,7),
    @(44,6,6),
    @(45,9,9),
    @(46,1,1),
    @(47,5,5),
    @(48,1,1),
    @(49,9,9),
    @(50,5,5),
    @(51,9,9),
    @(52,8,8),
    @(53,7,7),
    @(54,7,7),
    @(55,7,7),
    @(56,7,7),
    @(57,7,7),
    @(58,8,8),
    @(59,5,6),
    @(60,7,7),
    @(61,8,8),
    @(62,7,7),
    @(63,6,6),
    @(64,7,7),
    @(65,7,7),
    @(66,7,7),
    @(67,4,4),
    @(68,6,6),
    @(69,5,5),
    @(70,5,5),
    @(71,5,5),
    @(72,6,6),
    @(73,4,4),
    @(74,4,4),
    @(75,7,7),
    @(76,2,2)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
